$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update status of row 9 (Party Master in Transation IN form) from "In progress" to "Done"
$ws.Range("F9").Value = "Done"

# Add new column G with a note for row 9
$ws.Range("G9").Value = "pending in reports"

# Update selection to reflect the newly active / edited cell
$ws.Range("G9").Select()
